# Technology-cost-database.xlsx modifications
# - Remove the blank spacer rows in the "Database" sheet so the technology
#   rows are contiguous (rows shift up accordingly).
# - Replace the literal "?" placeholder in the Interest Rate column (F) of
#   the Photovoltaic Panels / Solar Collector / PV-thermal rows with the
#   numeric 5% interest rate (matching the other technologies), using the
#   same percentage number format already used elsewhere in column F.
# - Update the remembered cell selection to D19 (matches the new file).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Database")

# Remove the empty rows, bottom-up so earlier row numbers stay valid.
$ws.Rows.Item(16).Delete()
$ws.Rows.Item(14).Delete()
$ws.Rows.Item(11).Delete()
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(5).Delete()

# After the shift, the Photovoltaic Panels / Solar Collector / PV-thermal
# rows (previously 12, 13, 15, 17) now live at rows 8-11, with their
# Interest Rate cell (F) holding the text "?" instead of a real rate.
$ws.Range("F8:F11").Value = 0.05
$ws.Range("F8:F11").NumberFormat = "0%"

# Match the saved selection state recorded in the updated workbook.
$ws.Range("D19").Select()
